$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Gradient Boosting
$ws.Range("B2").Value = 0.8518993066906615
$ws.Range("C2").Value = 0.1690824380335535
$ws.Range("D2").Value = 0.7367906066536204
$ws.Range("E2").Value = 46051.7478635804

# Row 3 - Gradient Boosting (DeepWalk)
$ws.Range("B3").Value = 0.8550463317523451
$ws.Range("C3").Value = 0.1662750688711518
$ws.Range("D3").Value = 0.7446183953033269
$ws.Range("E3").Value = 45559.83866712334

# Row 4 - Linear Regression
$ws.Range("B4").Value = 0.6616725359431227
$ws.Range("C4").Value = 0.2938011756630212
$ws.Range("D4").Value = 0.474559686888454
$ws.Range("E4").Value = 69604.32066745973

# Row 5 - Linear Regression (DeepWalk)
$ws.Range("B5").Value = 0.6781529782285141
$ws.Range("C5").Value = 0.2919574826877662
$ws.Range("D5").Value = 0.4667318982387476
$ws.Range("E5").Value = 67887.89141201846

# Row 6 - Random Forest
$ws.Range("B6").Value = 0.8016464886531494
$ws.Range("C6").Value = 0.2124396871074345
$ws.Range("D6").Value = 0.6164383561643836
$ws.Range("E6").Value = 53295.12913925923

# Row 7 - Random Forest (DeepWalk)
$ws.Range("B7").Value = 0.816924994177087
$ws.Range("C7").Value = 0.2057637726716372
$ws.Range("D7").Value = 0.6457925636007827
$ws.Range("E7").Value = 51201.43120229485
